$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header E1 text ("Edge angle mean" -> "Mean edge angle")
$ws.Range("E1").Value = "Mean edge angle"

# 2. Copy existing row 14 (QSY_B_189 data) down to row 15, preserving styles/formats
$ws.Range("A14:H14").Copy($ws.Range("A15:H15"))

# 3. Overwrite row 14 with the new QSY_B_072 data
$ws.Range("A14").Value = "QSY_B_072"
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 0.62292363150254459
$ws.Range("D14").Value = 0.00015618284930595875
$ws.Range("E14").Value = 115.48421710219
$ws.Range("F14").Value = 8.0081252289049498
$ws.Range("G14").Value = 0.90212613413444998
$ws.Range("H14").Value = 1.02041757098209

# 4. Update the sheet selection to C17
$ws.Range("C17").Select()
